$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TimeRelay")

# Rename the "IsAvailable" header in D1 to "Availability"
$ws.Range("D1").Value = "Availability"

# The availability column now holds a plain number (0) instead of a boolean (TRUE)
$ws.Range("D2:D5").Value = 0

# Update the active selection shown when the sheet is next opened
$ws.Range("G10").Select()
